$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mojibake in the footnote text (shared string referenced by A103):
# corrupted accented characters ("í", "í", "ú", "ú") had been saved as ">" -
# restore the correct UTF-8 text.
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = `"Community of Sahel-Saharan States`";COMESA = `"Common Market for Eastern and Southern Africa`";EAC = `"East African Community`";ECCAS = `"Economic Community of Central African States`";ECOWAS = `"Economic Community of West African States`";IGAD = `"Intergovernmental Authority on Development`";SADC = `"Southern African Development Community`";UMA = `"Arab Maghreb Union`";PALOP = `"Países Africanos de Língua Oficial Portuguesa`";ASEAN = `"Association of Southeast Asian Nations`";MERCOSUR = `"Mercado Común del Sur`".EU27 = `"European Union (27 members)`".OECD = `"Organisation for Economic Co-operation and Development`"."

# Small recalculated-data correction for row 69 (Life expectancy at birth, 2020)
$ws.Range("C69").Value = 62.213457142857102

# Updated aggregate values for row 97 "Africa, Fragile States"
$ws.Range("C97").Value = 61.5719214285714
$ws.Range("D97").Value = 63.618450000000003
$ws.Range("E97").Value = 59.568453571428599
$ws.Range("F97").Value = 4.3557357142857098
$ws.Range("G97").Value = 47.112214285714302
$ws.Range("H97").Value = 67.454196428571507
$ws.Range("I97").Value = 163.34814285714299
$ws.Range("J97").Value = 25.819230769230799

# Updated aggregate values for row 98 "ROW, Fragile States"
$ws.Range("C98").Value = 69.959031249999995
$ws.Range("D98").Value = 72.6873875
$ws.Range("E98").Value = 67.407187500000006
$ws.Range("F98").Value = 2.6333625000000001
$ws.Range("G98").Value = 25.366475000000001
$ws.Range("H98").Value = 30.750556249999999
$ws.Range("I98").Value = 76.702743749999996
$ws.Range("J98").Value = 18.771428571428601
